$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J, matching existing header style (s="1")
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-63 for columns I and J
$data = @(
    @{ Row = 2; I = 5; J = 6 },
    @{ Row = 3; I = 9; J = 9 },
    @{ Row = 4; I = 12; J = 12 },
    @{ Row = 5; I = 7; J = 7 },
    @{ Row = 6; I = 6; J = 7 },
    @{ Row = 7; I = 7; J = 7 },
    @{ Row = 8; I = 5; J = 5 },
    @{ Row = 9; I = 7; J = 7 },
    @{ Row = 10; I = 9; J = 9 },
    @{ Row = 11; I = 7; J = 7 },
    @{ Row = 12; I = 7; J = 7 },
    @{ Row = 13; I = 7; J = 7 },
    @{ Row = 14; I = 7; J = 7 },
    @{ Row = 15; I = 8; J = 8 },
    @{ Row = 16; I = 8; J = 8 },
    @{ Row = 17; I = 7; J = 8 },
    @{ Row = 18; I = 9; J = 9 },
    @{ Row = 19; I = 8; J = 8 },
    @{ Row = 20; I = 8; J = 8 },
    @{ Row = 21; I = 8; J = 8 },
    @{ Row = 22; I = 6; J = 7 },
    @{ Row = 23; I = 9; J = 9 },
    @{ Row = 24; I = 6; J = 6 },
    @{ Row = 25; I = 7; J = 7 },
    @{ Row = 26; I = 5; J = 6 },
    @{ Row = 27; I = 4; J = 5 },
    @{ Row = 28; I = 8; J = 8 },
    @{ Row = 29; I = 9; J = 9 },
    @{ Row = 30; I = 9; J = 9 },
    @{ Row = 31; I = 9; J = 9 },
    @{ Row = 32; I = 8; J = 8 },
    @{ Row = 33; I = 8; J = 8 },
    @{ Row = 34; I = 7; J = 7 },
    @{ Row = 35; I = 7; J = 7 },
    @{ Row = 36; I = 7; J = 7 },
    @{ Row = 37; I = 8; J = 8 },
    @{ Row = 38; I = 7; J = 7 },
    @{ Row = 39; I = 7; J = 7 },
    @{ Row = 40; I = 7; J = 7 },
    @{ Row = 41; I = 6; J = 7 },
    @{ Row = 42; I = 6; J = 6 },
    @{ Row = 43; I = 8; J = 8 },
    @{ Row = 44; I = 7; J = 7 },
    @{ Row = 45; I = 8; J = 8 },
    @{ Row = 46; I = 7; J = 7 },
    @{ Row = 47; I = 7; J = 7 },
    @{ Row = 48; I = 7; J = 7 },
    @{ Row = 49; I = 4; J = 5 },
    @{ Row = 50; I = 7; J = 8 },
    @{ Row = 51; I = 8; J = 8 },
    @{ Row = 52; I = 5; J = 7 },
    @{ Row = 53; I = 5; J = 6 },
    @{ Row = 54; I = 7; J = 7 },
    @{ Row = 55; I = 9; J = 9 },
    @{ Row = 56; I = 7; J = 7 },
    @{ Row = 57; I = 7; J = 8 },
    @{ Row = 58; I = 5; J = 5 },
    @{ Row = 59; I = 8; J = 8 },
    @{ Row = 60; I = 5; J = 5 },
    @{ Row = 61; I = 7; J = 7 },
    @{ Row = 62; I = 3; J = 3 },
    @{ Row = 63; I = 2; J = 2 }
)

foreach ($row in $data) {
    $ws.Cells.Item($row.Row, 9).Value = $row.I
    $ws.Cells.Item($row.Row, 10).Value = $row.J
}
